# "Error Calculations and Plots"
# Re-roll the missing-data pattern on Sheet1:
#   - a handful of cells that were blanked out become populated again
#   - a handful of previously-populated cells become blanked out
#   - the "RM 232" and "SC 92" rows are removed entirely, shifting every
#     row below them up (dimension shrinks from A1:F35 to A1:F33)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level edits (original row numbers, applied before any rows move) ---

# Newly imputed values (previously blank -> numeric)
$ws.Range("D3").Value = -14.2
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("F12").Value = 17.45
$ws.Range("D21").Value = -14.3

# Newly missing values (previously numeric -> blank)
$ws.Range("F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("D23").ClearContents()

# "SC 132" (row 33) and "SC 193" (row 34) also get re-imputed values,
# ahead of the row deletions below shifting them up to rows 31 and 32.
$ws.Range("F33").Value = 17.18
$ws.Range("D34").Value = -14.7
$ws.Range("F34").Value = 17.39

# --- Remove the two rows that dropped out of this sample ("RM 232", "SC 92") ---
$ws.Rows(26).Delete()   # "RM 232" (row 27 "SC 92" shifts up to row 26)
$ws.Rows(27).Delete()   # "SC 92" (now at row 27)
